$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.975.22"
$ws.Range("E2").Value = "  -0.23%  "
$ws.Range("D3").Value = "2.623.90"
$ws.Range("E3").Value = "  -1.59%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'595.15"
$ws.Range("E5").Value = "  -0.38%  "
$ws.Range("D6").Value = "'166.82"
$ws.Range("E6").Value = "  +1.61%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("E8").Value = "  -2.50%  "
$ws.Range("E9").Value = "  -1.60%  "
$ws.Range("D10").Value = "'0.140"
$ws.Range("E10").Value = "  -1.52%  "
$ws.Range("E11").Value = "  +1.38%  "
$ws.Range("D12").Value = "'0.359"
$ws.Range("E12").Value = "  +0.27%  "
$ws.Range("E13").Value = "  -0.02%  "
$ws.Range("D14").Value = "'27.58"
$ws.Range("E14").Value = "  -0.80%  "
$ws.Range("D15").Value = "3.109.68"
$ws.Range("E16").Value = "  -1.12%  "
$ws.Range("D17").Value = "66.908.03"
$ws.Range("E17").Value = "  -0.44%  "
$ws.Range("D18").Value = "2.624.72"
$ws.Range("E18").Value = "  -1.33%  "
$ws.Range("E19").Value = "  +4.29%  "
$ws.Range("E20").Value = "  +6.93%  "
$ws.Range("D21").Value = "'355.83"
$ws.Range("E21").Value = "  -1.61%  "
$ws.Range("E22").Value = "  -1.89%  "
$ws.Range("E23").Value = "  -3.04%  "
$ws.Range("E24").Value = "  +8.48%  "
$ws.Range("E25").Value = "  -0.01%  "
$ws.Range("E26").Value = "  -5.80%  "
$ws.Range("D28").Value = "2.761.64"
$ws.Range("E28").Value = "  -1.45%  "
$ws.Range("E29").Value = "  +0.20%  "
$ws.Range("E30").Value = "  -1.48%  "
$ws.Range("D31").Value = "'547.84"
$ws.Range("E31").Value = "  -0.51%  "
$ws.Range("E32").Value = "  -0.80%  "
$ws.Range("E33").Value = "  -2.39%  "
$ws.Range("E35").Value = "  +5.10%  "
$ws.Range("E36").Value = "  +0.16%  "
$ws.Range("E37").Value = "  -5.02%  "
$ws.Range("E38").Value = "  +0.50%  "
$ws.Range("D39").Value = "'19.01"
$ws.Range("E40").Value = "  -1.82%  "
$ws.Range("E41").Value = "  -2.23%  "
$ws.Range("E42").Value = "  -2.26%  "
$ws.Range("E43").Value = "  +0.11%  "
$ws.Range("E45").Value = "  -0.11%  "
$ws.Range("E46").Value = "  -4.88%  "
$ws.Range("E47").Value = "  -0.82%  "
$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").Value = "'0.577"
$ws.Range("E48").Value = "  -1.30%  "
$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").Value = "'151.23"
$ws.Range("E49").Value = "  -1.29%  "
$ws.Range("E50").Value = "  -1.70%  "
$ws.Range("E51").Value = "  -0.92%  "

# Restore default style on cells forced to text via quote-prefix, to avoid lingering quotePrefix style flag
$textForcedCells = @("D5","D6","D10","D12","D14","D21","D31","D39","D48","D49")
foreach ($addr in $textForcedCells) {
    $ws.Range($addr).Style = "Normal"
}
